$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 44

$ws.Cells.Item($row, 1).Value = 44357.77969643586
$ws.Cells.Item($row, 2).Value = 76895
$ws.Cells.Item($row, 3).Value = 64695
$ws.Cells.Item($row, 4).Value = 3392
$ws.Cells.Item($row, 5).Value = 2120
$ws.Cells.Item($row, 6).Value = 1482
$ws.Cells.Item($row, 7).Value = 20294
$ws.Cells.Item($row, 8).Value = 1398
$ws.Cells.Item($row, 9).Value = 880
$ws.Cells.Item($row, 10).Value = 186

# Apply the same style as the other date cells in column A (style index 2 -> numFmt "yyyy-mm-dd HH:mm:ss UTC")
$ws.Cells.Item($row, 1).NumberFormat = $ws.Cells.Item($row - 1, 1).NumberFormat
